$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) for rows 2-45 ---
$ws.Range("D2").Value = "29.216.81"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.855.48"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.04"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6989"
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9997"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07720"
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3076"
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.78"
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07806"
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("D12").Value = "1.865.08"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "92.25"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.093"
$ws.Range("E14").Value = "  -1.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6867"
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.517"
$ws.Range("E16").Value = "  +2.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008415"
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("D18").Value = "29.233.74"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "249.34"
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("D20").Value = "2.113.81"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.80"
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9994"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.506"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9996"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1500"
$ws.Range("E25").Value = "  -3.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.65"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.840"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.48"
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.558"
$ws.Range("E29").Value = "  +4.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.234"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.176"
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.193"
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05202"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7600"
$ws.Range("E34").Value = "  +2.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.842"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.165"
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.709"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01862"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").Value = "1.216.95"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.719"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8971"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.72"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9987"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.528"
$ws.Range("E44").Value = "  -12.18%  "
$ws.Range("D45").Value = "2.012.64"
$ws.Range("E45").Value = "  -2.00%  "

# --- Rows 46-49 shifted (re-ranked coins) with new Coin/Link/Price/Volume ---
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000124"
$ws.Range("E46").Value = "  -3.32%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.19"
$ws.Range("E47").Value = "  -8.87%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5177"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.534"
$ws.Range("E49").Value = "  +1.58%  "

# --- Rows 50-51: only Volume(1h) changes ---
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("E51").Value = "  +0.55%  "
